# Apply the CodeSystem-HistoGradeCS metadata refresh:
#  - Insert a new "Jurisdiction" property row (with an empty value) right
#    after the existing "Contact" row on the Metadata sheet, pushing every
#    row below it down by one (Count ends up on row 22 instead of 21).
#  - Refresh the "Date" property's timestamp value.
#
# The "Concepts" sheet is left untouched; its shared-string references
# shift automatically as a side effect of the Metadata sheet's edits.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$lastRow  = 21   # last populated row before the edit
$insertAt = 11   # new "Jurisdiction" row goes here (right after "Contact")

# --- 1. Note which rows have a blank column-B up front -------------------
#        (pasting a blank source cell is a no-op in this engine, so those
#        destinations need an explicit ClearContents instead).
$blankB = @{}
for ($r = $insertAt; $r -le $lastRow; $r++) {
    $v = $ws.Cells.Item($r, 2).Value()
    $blankB[$r] = ($v -eq $null) -or ($v -eq "")
}

# --- 2. Push rows $insertAt..$lastRow down by one, bottom-up so the ------
#        source row is never clobbered before it has been read/copied.
for ($r = $lastRow; $r -ge $insertAt; $r--) {
    $destRow  = $r + 1
    $srcRange = $ws.Range("A" + $r + ":B" + $r)
    $dstRange = $ws.Range("A" + $destRow + ":B" + $destRow)

    # Formats first (also covers brand-new row 22, which has no style yet).
    $srcRange.Copy()
    $dstRange.PasteSpecial(-4122)   # xlPasteFormats
    $excel.CutCopyMode = 0

    # Column A is always non-blank text; paste-special keeps it text-typed.
    $srcRange.Copy()
    $ws.Range("A" + $destRow).PasteSpecial(-4163)   # xlPasteValues
    $excel.CutCopyMode = 0

    # Column B: value-paste only when genuinely non-blank, otherwise clear.
    if ($blankB[$r]) {
        $ws.Cells.Item($destRow, 2).ClearContents()
    } else {
        $ws.Range("B" + $r).Copy()
        $ws.Range("B" + $destRow).PasteSpecial(-4163)
        $excel.CutCopyMode = 0
    }
}

# --- 3. Write the new "Jurisdiction" row into the freed-up row 11 --------
$ws.Cells.Item($insertAt, 1).Value = "Jurisdiction"
$ws.Cells.Item($insertAt, 2).Value = ""

# --- 4. Refresh the "Date" property value --------------------------------
$ws.Cells.Item(8, 2).Value = "2024-09-17T19:55:11+00:00"
